$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Traducciones")

# Fix typo: hour_rage -> hour_range
$ws.Cells.Item(133, 1).Value = "hour_range"

# Row 135: select_student
$ws.Cells.Item(135, 1).Value = "select_student"
$ws.Cells.Item(135, 2).Value = "Select students"
$ws.Cells.Item(135, 3).Value = "Selecciona un estudiante"
$ws.Cells.Item(135, 4).Value = "Aukeratu ikasle bat"
$ws.Cells.Item(135, 5).Value = "Selecciona un estudiant"
$ws.Cells.Item(135, 6).Value = "Selecciona un alumno"

# Row 136: successful_password_restart
$ws.Cells.Item(136, 1).Value = "successful_password_restart"
$ws.Cells.Item(136, 2).Value = "Student password has been reset"
$ws.Cells.Item(136, 3).Value = "La contraseña del estudiante ha sido reestablecida"
$ws.Cells.Item(136, 4).Value = "Ikaslearen pasahitza berrezarri da"
$ws.Cells.Item(136, 5).Value = "La contrasenya de l'estudiant ha estat restablida"
$ws.Cells.Item(136, 6).Value = "Restableceuse o contrasinal do alumno"

# Row 137: error_restarting_student_password
$ws.Cells.Item(137, 1).Value = "error_restarting_student_password"
$ws.Cells.Item(137, 2).Value = "Error resetting student password"
$ws.Cells.Item(137, 3).Value = "Error al restablecer la contraseña del estudiante"
$ws.Cells.Item(137, 4).Value = "Errore bat gertatu da ikaslearen pasahitza berrezartzean"
$ws.Cells.Item(137, 5).Value = "Error en restablir la contrasenya de l'estudiant"
$ws.Cells.Item(137, 6).Value = "Produciuse un erro ao restablecer o contrasinal do alumno"

# Row 138: student_grades
$ws.Cells.Item(138, 1).Value = "student_grades"
$ws.Cells.Item(138, 2).Value = "Student grades"
$ws.Cells.Item(138, 3).Value = "Notas estudiantes"
$ws.Cells.Item(138, 4).Value = "Ikasleen kalifikazioak"
$ws.Cells.Item(138, 5).Value = "Qualificacions estudiants"
$ws.Cells.Item(138, 6).Value = "Cualificacións dos estudantes"

# Row 139: exercises_get_fail
$ws.Cells.Item(139, 1).Value = "exercises_get_fail"
$ws.Cells.Item(139, 2).Value = "Error getting exercises"
$ws.Cells.Item(139, 3).Value = "Error al obtener ejercicios"
$ws.Cells.Item(139, 4).Value = "Errore bat gertatu da ariketak eskuratzean"
$ws.Cells.Item(139, 5).Value = "Error en obtenir exercicis"
$ws.Cells.Item(139, 6).Value = "Erro ao conseguir exercicios"

# Row 140: clear A140 and B140 (become fully blank cells, no style info retained as value cells)
$ws.Cells.Item(140, 1).ClearContents()
$ws.Cells.Item(140, 2).ClearContents()
